# LivInParis "Peuplement" workbook update:
#   - Account sheet: rename the "email" column to "user_name" and replace
#     each "userN@example.com" value with the bare username "userN".
#   - Make the Account sheet the active/selected tab (cell A2 selected),
#     instead of the Chef sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account")

# Header: email -> user_name (password column B is untouched)
$ws.Cells.Item(1, 1).Value = "user_name"

# Rows 2..51 hold user0..user49; strip the "@example.com" suffix from column A.
for ($i = 0; $i -le 49; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "user$i"
}

# Switch the active sheet/selection from Chef to Account.
$ws.Activate()
$ws.Range("A2").Select()
